# Remove the "NA" placeholder values from the Topic column (D) of the
# media coverage sheet. Rows where the Topic was genuinely known keep
# their value; rows whose Topic cell literally equals "NA" have that
# cell cleared entirely (matching the source data export after the
# author reran the topic-classification / distribution step).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq "NA") {
        $cell.ClearContents()
    }
}
